$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.004.44"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "3.502.35"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.92%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +3.77%  "

$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "4.103.85"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.69%  "

$ws.Range("D15").Value = "67.020.28"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "3.483.60"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "

$ws.Range("E30").Value = "  -2.23%  "

$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.880"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +2.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("D42").Value = "2.835.65"
$ws.Range("E42").Value = "  +2.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("E45").Value = "  +2.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0303"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "339.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.844"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.47%  "

$ws.Range("E51").Value = "  -0.76%  "
